# Auto-generated edit script: update market-price derived columns (H-N)
# per scheduled market-data refresh, across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 444.125
$ws.Cells.Item(6, 9).Value = 444.125
$ws.Cells.Item(6, 11).Value = 1332.375
$ws.Cells.Item(6, 13).Value = -1220.375
$ws.Cells.Item(28, 8).Value = 830.3333
$ws.Cells.Item(28, 9).Value = 871.75
$ws.Cells.Item(28, 10).Value = 747.5
$ws.Cells.Item(28, 11).Value = 871.75
$ws.Cells.Item(28, 12).Value = 747.5
$ws.Cells.Item(28, 13).Value = -386.75
$ws.Cells.Item(28, 14).Value = -1717.5
$ws.Cells.Item(33, 8).Value = 762.64703
$ws.Cells.Item(33, 9).Value = 833.26666
$ws.Cells.Item(33, 10).Value = 233
$ws.Cells.Item(33, 11).Value = 833.26666
$ws.Cells.Item(33, 12).Value = 233
$ws.Cells.Item(33, 13).Value = -604.26666
$ws.Cells.Item(33, 14).Value = -691
$ws.Cells.Item(42, 8).Value = 156.18182
$ws.Cells.Item(42, 9).Value = 98.666664
$ws.Cells.Item(42, 10).Value = 415
$ws.Cells.Item(42, 11).Value = 295.999992
$ws.Cells.Item(42, 12).Value = 1245
$ws.Cells.Item(42, 13).Value = -65.99999200000002
$ws.Cells.Item(42, 14).Value = -1705

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17392.893
$ws.Cells.Item(32, 9).Value = 4998.019
$ws.Cells.Item(32, 10).Value = 33816.1
$ws.Cells.Item(32, 11).Value = 4998.019
$ws.Cells.Item(32, 12).Value = 33816.1
$ws.Cells.Item(32, 13).Value = -4711.019
$ws.Cells.Item(32, 14).Value = -34390.1
$ws.Cells.Item(35, 8).Value = 1007.4
$ws.Cells.Item(35, 9).Value = 1007.4
$ws.Cells.Item(35, 11).Value = 1007.4
$ws.Cells.Item(35, 13).Value = -601.4
$ws.Cells.Item(44, 8).Value = 19800
$ws.Cells.Item(44, 10).Value = 19800
$ws.Cells.Item(44, 12).Value = 19800
$ws.Cells.Item(44, 14).Value = -20776
$ws.Cells.Item(74, 8).Value = 4379.05
$ws.Cells.Item(74, 9).Value = 5703.4443
$ws.Cells.Item(74, 11).Value = 5703.4443
$ws.Cells.Item(74, 13).Value = -4829.4443
$ws.Cells.Item(77, 8).Value = 4379.05
$ws.Cells.Item(77, 9).Value = 5703.4443
$ws.Cells.Item(77, 11).Value = 28517.2215
$ws.Cells.Item(77, 13).Value = -24149.2215
$ws.Cells.Item(132, 8).Value = 1878.5349
$ws.Cells.Item(132, 9).Value = 1617.2941
$ws.Cells.Item(132, 10).Value = 2865.4443
$ws.Cells.Item(132, 11).Value = 4851.8823
$ws.Cells.Item(132, 12).Value = 8596.332900000001
$ws.Cells.Item(132, 13).Value = -2321.8823
$ws.Cells.Item(132, 14).Value = -13656.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 4700
$ws.Cells.Item(22, 9).Value = 4550
$ws.Cells.Item(22, 10).Value = 5000
$ws.Cells.Item(22, 11).Value = 4550
$ws.Cells.Item(22, 12).Value = 5000
$ws.Cells.Item(22, 13).Value = -4377
$ws.Cells.Item(22, 14).Value = -5346

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3775885.2
$ws.Cells.Item(31, 9).Value = 10001156
$ws.Cells.Item(31, 10).Value = 2993.7576
$ws.Cells.Item(31, 11).Value = 10001156
$ws.Cells.Item(31, 12).Value = 2993.7576
$ws.Cells.Item(31, 13).Value = -10000861
$ws.Cells.Item(31, 14).Value = -3583.7576
$ws.Cells.Item(34, 8).Value = 3775885.2
$ws.Cells.Item(34, 9).Value = 10001156
$ws.Cells.Item(34, 10).Value = 2993.7576
$ws.Cells.Item(34, 11).Value = 10001156
$ws.Cells.Item(34, 12).Value = 2993.7576
$ws.Cells.Item(34, 13).Value = -10000954
$ws.Cells.Item(34, 14).Value = -3397.7576

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 71
$ws.Cells.Item(13, 9).Value = 71
$ws.Cells.Item(13, 11).Value = 213
$ws.Cells.Item(13, 13).Value = -45
$ws.Cells.Item(14, 8).Value = 112.57143
$ws.Cells.Item(14, 9).Value = 112.57143
$ws.Cells.Item(14, 11).Value = 337.71429
$ws.Cells.Item(14, 13).Value = -164.71429
$ws.Cells.Item(17, 8).Value = 1249.1111
$ws.Cells.Item(17, 9).Value = 150
$ws.Cells.Item(17, 10).Value = 1563.1428
$ws.Cells.Item(17, 11).Value = 450
$ws.Cells.Item(17, 12).Value = 4689.428400000001
$ws.Cells.Item(17, 13).Value = -281
$ws.Cells.Item(17, 14).Value = -5027.428400000001
$ws.Cells.Item(113, 8).Value = 754.5854
$ws.Cells.Item(113, 9).Value = 415.08
$ws.Cells.Item(113, 10).Value = 1285.0625
$ws.Cells.Item(113, 11).Value = 1245.24
$ws.Cells.Item(113, 12).Value = 3855.1875
$ws.Cells.Item(113, 13).Value = 924.76
$ws.Cells.Item(113, 14).Value = -8195.1875
$ws.Cells.Item(122, 8).Value = 716.8261
$ws.Cells.Item(122, 9).Value = 363.3846
$ws.Cells.Item(122, 10).Value = 1176.3
$ws.Cells.Item(122, 11).Value = 3270.4614
$ws.Cells.Item(122, 12).Value = 10586.7
$ws.Cells.Item(122, 13).Value = -820.4613999999997
$ws.Cells.Item(122, 14).Value = -15486.7
$ws.Cells.Item(140, 8).Value = 1522.3
$ws.Cells.Item(140, 10).Value = 4500
$ws.Cells.Item(140, 12).Value = 13500
$ws.Cells.Item(140, 14).Value = -23860

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 24402
$ws.Cells.Item(12, 9).Value = 24402
$ws.Cells.Item(12, 11).Value = 24402
$ws.Cells.Item(12, 13).Value = -24262

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 1400080
$ws.Cells.Item(2, 9).Value = 100
$ws.Cells.Item(2, 10).Value = 1750075
$ws.Cells.Item(2, 11).Value = 100
$ws.Cells.Item(2, 12).Value = 1750075
$ws.Cells.Item(2, 13).Value = 12
$ws.Cells.Item(2, 14).Value = -1750299
$ws.Cells.Item(3, 8).Value = 2995
$ws.Cells.Item(3, 10).Value = 2995
$ws.Cells.Item(3, 12).Value = 2995
$ws.Cells.Item(3, 14).Value = -3219
$ws.Cells.Item(15, 8).Value = 2995
$ws.Cells.Item(15, 10).Value = 2995
$ws.Cells.Item(15, 12).Value = 2995
$ws.Cells.Item(15, 14).Value = -3335
$ws.Cells.Item(53, 8).Value = 11000
$ws.Cells.Item(53, 10).Value = 11000
$ws.Cells.Item(53, 12).Value = 11000
$ws.Cells.Item(53, 14).Value = -12036
$ws.Cells.Item(132, 8).Value = 14714715
$ws.Cells.Item(132, 9).Value = 31265438
$ws.Cells.Item(132, 11).Value = 93796314
$ws.Cells.Item(132, 13).Value = -93793784

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(13, 8).Value = 533
$ws.Cells.Item(13, 9).Value = 549.5
$ws.Cells.Item(13, 10).Value = 500
$ws.Cells.Item(13, 11).Value = 549.5
$ws.Cells.Item(13, 12).Value = 500
$ws.Cells.Item(13, 13).Value = -409.5
$ws.Cells.Item(13, 14).Value = -780
$ws.Cells.Item(15, 8).Value = 7864.143
$ws.Cells.Item(15, 10).Value = 7864.143
$ws.Cells.Item(15, 12).Value = 7864.143
$ws.Cells.Item(15, 14).Value = -8440.143
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(126, 8).Value = 5575.5
$ws.Cells.Item(126, 9).Value = 7217.3335
$ws.Cells.Item(126, 11).Value = 21652.0005
$ws.Cells.Item(126, 13).Value = -19182.0005
$ws.Cells.Item(132, 8).Value = 2152.75
$ws.Cells.Item(132, 9).Value = 1040.8
$ws.Cells.Item(132, 10).Value = 2445.3684
$ws.Cells.Item(132, 11).Value = 3122.4
$ws.Cells.Item(132, 12).Value = 7336.1052
$ws.Cells.Item(132, 13).Value = -592.3999999999996
$ws.Cells.Item(132, 14).Value = -12396.1052
$ws.Cells.Item(136, 8).Value = 2026.2727
$ws.Cells.Item(136, 9).Value = 1430.4
$ws.Cells.Item(136, 10).Value = 2522.8333
$ws.Cells.Item(136, 11).Value = 3204
$ws.Cells.Item(136, 12).Value = 7568.499899999999
$ws.Cells.Item(136, 13).Value = -1741.200000000001
$ws.Cells.Item(136, 14).Value = -12668.4999

# N53 on WVR is removed entirely in the target state (not merely zeroed).
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(53, 14).ClearContents()

